$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.461.02"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.089.29"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.68"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.79"
$ws.Range("E6").Value = "  +3.16%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("E8").Value = "  +3.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.075.12"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.84"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000241"
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.35"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.600.84"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.15"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.503.54"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.081.00"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.86"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.64"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.717"
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.10"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.72"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.67"
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.19"
$ws.Range("E32").Value = "  -3.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.113"
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.37"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0845"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("E37").Value = "  +2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.09"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("E39").Value = "  -5.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.34"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.55"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "441.76"
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.816.43"
$ws.Range("E45").Value = "  -4.20%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.110"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.98"
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.04"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.09"
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("E51").Value = "  +0.31%  "
